$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range first so numeric-looking strings
# (e.g. "587.47") are stored as text, matching the original inlineStr cells,
# instead of being auto-converted to numbers by Excel.
$dataRange = $ws.Range("B2:E51")
$dataRange.NumberFormat = "@"

$ws.Range('D2').Value = '69.302.14'
$ws.Range('E2').Value = '  +2.18%  '
$ws.Range('D3').Value = '3.388.05'
$ws.Range('E3').Value = '  +1.62%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '587.47'
$ws.Range('E5').Value = '  +0.68%  '
$ws.Range('D6').Value = '180.75'
$ws.Range('E6').Value = '  +2.63%  '
$ws.Range('E7').Value = '  -0.02%  '
$ws.Range('E8').Value = '  +0.70%  '
$ws.Range('D9').Value = '0.197'
$ws.Range('E9').Value = '  +8.83%  '
$ws.Range('E10').Value = '  +1.41%  '
$ws.Range('D11').Value = '48.79'
$ws.Range('E11').Value = '  +4.81%  '
$ws.Range('E12').Value = '  +4.64%  '
$ws.Range('D13').Value = '684.20'
$ws.Range('E13').Value = '  -0.93%  '
$ws.Range('D14').Value = '8.64'
$ws.Range('D15').Value = '3.932.34'
$ws.Range('E15').Value = '  +1.55%  '
$ws.Range('D16').Value = '69.331.66'
$ws.Range('E16').Value = '  +2.23%  '
$ws.Range('B17').Value = 'WrappedEther'
$ws.Range('C17').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D17').Value = '3.390.66'
$ws.Range('E17').Value = '  +1.95%  '
$ws.Range('B18').Value = 'TRON'
$ws.Range('C18').Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range('D18').Value = '0.120'
$ws.Range('E18').Value = '  +1.77%  '
$ws.Range('D19').Value = '17.71'
$ws.Range('E19').Value = '  +1.73%  '
$ws.Range('D20').Value = '11.38'
$ws.Range('E20').Value = '  +2.95%  '
$ws.Range('D21').Value = '0.901'
$ws.Range('E21').Value = '  +0.62%  '
$ws.Range('D22').Value = '5.44'
$ws.Range('E22').Value = '  +1.35%  '
$ws.Range('E23').Value = '  +0.36%  '
$ws.Range('D24').Value = '103.70'
$ws.Range('E24').Value = '  +5.19%  '
$ws.Range('E25').Value = '  +1.36%  '
$ws.Range('E26').Value = '  +1.43%  '
$ws.Range('D27').Value = '9.61'
$ws.Range('E27').Value = '  +1.35%  '
$ws.Range('D28').Value = '34.25'
$ws.Range('E28').Value = '  +3.83%  '
$ws.Range('D29').Value = '8.73'
$ws.Range('E29').Value = '  +1.90%  '
$ws.Range('E30').Value = '  -1.65%  '
$ws.Range('D31').Value = '3.69'
$ws.Range('E31').Value = '  +10.40%  '
$ws.Range('E32').Value = '  +1.77%  '
$ws.Range('D33').Value = '555.54'
$ws.Range('E33').Value = '  -2.91%  '
$ws.Range('E34').Value = '  +0.89%  '
$ws.Range('D35').Value = '58.13'
$ws.Range('E35').Value = '  +1.34%  '
$ws.Range('E36').Value = '  +0.05%  '
$ws.Range('D37').Value = '3.698.63'
$ws.Range('E37').Value = '  -0.46%  '
$ws.Range('E38').Value = '  +7.54%  '
$ws.Range('D39').Value = '35.03'
$ws.Range('E39').Value = '  +3.17%  '
$ws.Range('D40').Value = '3.24'
$ws.Range('E40').Value = '  +1.65%  '
$ws.Range('B41').Value = 'PEPE'
$ws.Range('C41').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D41').Value = '0.0₃0704'
$ws.Range('E41').Value = '  +4.19%  '
$ws.Range('B42').Value = 'Fetch.AI'
$ws.Range('C42').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D42').Value = '2.68'
$ws.Range('E42').Value = '  +1.52%  '
$ws.Range('E43').Value = '  +1.01%  '
$ws.Range('D44').Value = '0.0423'
$ws.Range('E44').Value = '  +4.26%  '
$ws.Range('D45').Value = '3.25'
$ws.Range('E45').Value = '  -2.50%  '
$ws.Range('E46').Value = '  +0.11%  '
$ws.Range('E47').Value = '  +0.83%  '
$ws.Range('B48').Value = 'FirstDigitalUSD'
$ws.Range('C48').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D48').Value = '1.00'
$ws.Range('E48').Value = '  +0.00%  '
$ws.Range('B49').Value = 'Mantle'
$ws.Range('C49').Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range('D49').Value = '1.38'
$ws.Range('E49').Value = '  +4.84%  '
$ws.Range('D50').Value = '132.49'
$ws.Range('E50').Value = '  +2.13%  '
$ws.Range('B51').Value = 'THORChain'
$ws.Range('C51').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D51').Value = '7.51'
$ws.Range('E51').Value = '  +1.40%  '

# Restore the default cell style so no stray number-format styling remains
$dataRange.Style = "Normal"
